$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 data: game 13's entry
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Ruby Soho"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = "Wise"
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 6
$ws.Range("I15").Value = "Sell 3"
$ws.Range("J15").Value = "Tax+4:Prussia;Sail+3:Windy;Advisors:Patmos"

# Update selection to match the new active cell
$ws.Range("J15").Select()
